# Updates cryptos list price/volume columns (generated from upstream data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.636.51'
$ws.Range("E2").Value = '  -1.22%  '
$ws.Range("D3").Value = '1.849.60'
$ws.Range("E3").Value = '  -0.90%  '
$ws.Range("E4").Value = '  -0.52%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.79'
$ws.Range("E5").Value = '  -1.52%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9994'
$ws.Range("E6").Value = '  -0.51%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4244'
$ws.Range("E7").Value = '  -1.71%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3642'
$ws.Range("E8").Value = '  -1.66%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '44.46'
$ws.Range("E9").Value = '  -0.07%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07307'
$ws.Range("E10").Value = '  -1.16%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8781'
$ws.Range("E11").Value = '  -5.36%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '20.78'
$ws.Range("E12").Value = '  -1.85%  '
$ws.Range("D13").Value = '1.838.51'
$ws.Range("E13").Value = '  -7.83%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.344'
$ws.Range("E14").Value = '  -1.22%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.534'
$ws.Range("E15").Value = '  -2.64%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.06899'
$ws.Range("E16").Value = '  +0.66%  '
$ws.Range("E17").Value = '  -0.51%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '79.30'
$ws.Range("E18").Value = '  -0.98%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000008904'
$ws.Range("E19").Value = '  -0.94%  '
$ws.Range("E20").Value = '  -0.30%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.39'
$ws.Range("E21").Value = '  -1.99%  '
$ws.Range("D22").Value = '27.679.87'
$ws.Range("E22").Value = '  -1.09%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.990'
$ws.Range("E23").Value = '  -2.15%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.45'
$ws.Range("E24").Value = '  -4.67%  '
$ws.Range("D25").Value = '2.091.88'
$ws.Range("E25").Value = '  -4.90%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.980'
$ws.Range("E26").Value = '  -3.35%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '152.72'
$ws.Range("E27").Value = '  -0.76%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.96'
$ws.Range("E28").Value = '  +2.68%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '121.27'
$ws.Range("E29").Value = '  +7.41%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.277'
$ws.Range("E30").Value = '  -3.62%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.888'
$ws.Range("E31").Value = '  +12.12%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08887'
$ws.Range("E32").Value = '  -0.90%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7685'
$ws.Range("E33").Value = '  -4.20%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.570'
$ws.Range("E34").Value = '  -3.86%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.966'
$ws.Range("E35").Value = '  +0.48%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.105'
$ws.Range("E36").Value = '  -5.43%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9983'
$ws.Range("E37").Value = '  -0.73%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.095'
$ws.Range("E38").Value = '  -2.09%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05366'
$ws.Range("E39").Value = '  -2.14%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01937'
$ws.Range("E40").Value = '  -1.54%  '
$ws.Range("E41").Value = '  -6.28%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5122'
$ws.Range("E42").Value = '  -1.71%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.901'
$ws.Range("E43").Value = '  -1.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1652'
$ws.Range("E44").Value = '  -1.58%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.310'
$ws.Range("E45").Value = '  -4.52%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.06541'
$ws.Range("E46").Value = '  -2.41%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4762'
$ws.Range("E47").Value = '  -1.94%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '10.36'
$ws.Range("E48").Value = '  -1.15%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '104.94'
$ws.Range("E49").Value = '  -1.36%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.9984'
$ws.Range("E50").Value = '  -0.59%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.629'
$ws.Range("E51").Value = '  -2.07%  '
